$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their text representation (e.g. trailing zeros, grouped
# "thousand.thousand" style price strings) instead of being auto-converted
# to numbers by Excel when the new value looks numeric.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.224.49"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.571.20"
$ws.Range("E3").Value = "  +0.63%  "
$ws.Range("E4").Value = "  +0.60%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.63"
$ws.Range("E5").Value = "  +2.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.493"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("E7").Value = "  +0.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.14"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.250"
$ws.Range("E9").Value = "  +0.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0601"
$ws.Range("E10").Value = "  +0.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0869"
$ws.Range("E11").Value = "  +1.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.793.65"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.582.31"
$ws.Range("E13").Value = "  +1.71%  "
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "62.43"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.193.12"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0₃0703"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "216.48"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.42"
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("E22").Value = "  +1.14%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.07"
$ws.Range("E25").Value = "  +0.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.73"
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.10"
$ws.Range("E27").Value = "  +0.29%  "
$ws.Range("E28").Value = "  +2.17%  "
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("E30").Value = "  +2.44%  "
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.452.83"
$ws.Range("E33").Value = "  +1.99%  "
$ws.Range("E34").Value = "  +1.88%  "
$ws.Range("E35").Value = "  +5.05%  "
$ws.Range("E36").Value = "  +0.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.36"
$ws.Range("E37").Value = "  +1.46%  "
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.536"
$ws.Range("E39").Value = "  +0.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.83"
$ws.Range("E40").Value = "  +2.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.810"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.76"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.706.71"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.96"
$ws.Range("E48").Value = "  -1.79%  "
$ws.Range("E49").Value = "  +3.95%  "
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0959"
$ws.Range("E51").Value = "  +0.39%  "
